$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 11 (Leve Item ID 5533)
$ws.Range("H11").Value = 1.3333334
$ws.Range("I11").Value = 1.3333334
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 1.3333334
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = 138.6666666

# Row 19 (Leve Item ID 7015)
$ws.Range("H19").Value = 1524.75
$ws.Range("I19").Value = 1466.6666
$ws.Range("J19").Value = 1699
$ws.Range("K19").Value = 1466.6666
$ws.Range("L19").Value = 1699
$ws.Range("M19").Value = -1291.6666
$ws.Range("N19").Value = -2049

# Row 32 (Leve Item ID 5484)
$ws.Range("H32").Value = 1499
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 1499
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 1499
$ws.Range("M32").Value = ""
$ws.Range("N32").Value = -2151

# Row 33 (Leve Item ID 5512)
$ws.Range("H33").Value = 329
$ws.Range("I33").Value = 329
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 329
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -100
$ws.Range("N33").Value = ""

# Row 70 (Leve Item ID 12604)
$ws.Range("H70").Value = 2859.6667
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 2859.6667
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 8579.000100000001
$ws.Range("M70").Value = ""
$ws.Range("N70").Value = -9119.000100000001

# Row 73 (Leve Item ID 12604)
$ws.Range("H73").Value = 2859.6667
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 2859.6667
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 8579.000100000001
$ws.Range("M73").Value = ""
$ws.Range("N73").Value = -10451.0001

# Row 80 (Leve Item ID 12605)
$ws.Range("H80").Value = 1184
$ws.Range("I80").Value = 333
$ws.Range("J80").Value = 1467.6666
$ws.Range("K80").Value = 999
$ws.Range("L80").Value = 4402.9998
$ws.Range("M80").Value = -1
$ws.Range("N80").Value = -6398.9998

# Row 83 (Leve Item ID 12605)
$ws.Range("H83").Value = 1184
$ws.Range("I83").Value = 333
$ws.Range("J83").Value = 1467.6666
$ws.Range("K83").Value = 2997
$ws.Range("L83").Value = 13208.9994
$ws.Range("M83").Value = 1995
$ws.Range("N83").Value = -23192.9994

# Row 96 (Leve Item ID 19894)
$ws.Range("H96").Value = 584
$ws.Range("I96").Value = 391
$ws.Range("J96").Value = 777
$ws.Range("K96").Value = 1173
$ws.Range("L96").Value = 2331
$ws.Range("M96").Value = 200
$ws.Range("N96").Value = -5077

# Row 100 (Leve Item ID 19906)
$ws.Range("H100").Value = 3671.6667
$ws.Range("I100").Value = 2149.6667
$ws.Range("J100").Value = 8237.666999999999
$ws.Range("K100").Value = 2149.6667
$ws.Range("L100").Value = 8237.666999999999
$ws.Range("M100").Value = -1608.6667
$ws.Range("N100").Value = -9319.666999999999

# Row 132 (Leve Item ID 44049)
$ws.Range("H132").Value = 949.3889
$ws.Range("I132").Value = 484.64517
$ws.Range("J132").Value = 3830.8
$ws.Range("K132").Value = 1453.93551
$ws.Range("L132").Value = 11492.4
$ws.Range("M132").Value = 1076.06449
$ws.Range("N132").Value = -16552.4

# Row 135 (Leve Item ID 44047)
$ws.Range("H135").Value = 20834014
$ws.Range("I135").Value = 21739836
$ws.Range("J135").Value = 81
$ws.Range("K135").Value = 195658524
$ws.Range("L135").Value = 729
$ws.Range("M135").Value = -195655989
$ws.Range("N135").Value = -5799

$ws = $wb.Worksheets.Item("ARM")
# Row 2 (Leve Item ID 27713)
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = ""
$ws.Range("N2").Value = ""

# Row 74 (Leve Item ID 44000)
$ws.Range("H74").Value = 760.9636
$ws.Range("I74").Value = 615.7143
$ws.Range("J74").Value = 1947.1666
$ws.Range("K74").Value = 615.7143
$ws.Range("L74").Value = 1947.1666
$ws.Range("M74").Value = 258.2857
$ws.Range("N74").Value = -3695.1666

# Row 77 (Leve Item ID 44000)
$ws.Range("H77").Value = 760.9636
$ws.Range("I77").Value = 615.7143
$ws.Range("J77").Value = 1947.1666
$ws.Range("K77").Value = 3078.5715
$ws.Range("L77").Value = 9735.833000000001
$ws.Range("M77").Value = 1289.4285
$ws.Range("N77").Value = -18471.833

# Row 116 (Leve Item ID 27713)
$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = ""
$ws.Range("N116").Value = ""

# Row 132 (Leve Item ID 43997)
$ws.Range("H132").Value = 19233126
$ws.Range("I132").Value = 55556508
$ws.Range("J132").Value = 3101
$ws.Range("K132").Value = 166669524
$ws.Range("L132").Value = 9303
$ws.Range("M132").Value = -166666994
$ws.Range("N132").Value = -14363

$ws = $wb.Worksheets.Item("BSM")
# Row 3 (Leve Item ID 27713)
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = ""
$ws.Range("N3").Value = ""

# Row 80 (Leve Item ID 13747)
$ws.Range("H80").Value = 294.66666
$ws.Range("I80").Value = 421.5
$ws.Range("J80").Value = 41
$ws.Range("K80").Value = 421.5
$ws.Range("L80").Value = 41
$ws.Range("M80").Value = 576.5
$ws.Range("N80").Value = -2037

# Row 83 (Leve Item ID 13747)
$ws.Range("H83").Value = 294.66666
$ws.Range("I83").Value = 421.5
$ws.Range("J83").Value = 41
$ws.Range("K83").Value = 2107.5
$ws.Range("L83").Value = 205
$ws.Range("M83").Value = 2884.5
$ws.Range("N83").Value = -10189

# Row 86 (Leve Item ID 12526)
$ws.Range("H86").Value = 2726.5
$ws.Range("I86").Value = 703
$ws.Range("J86").Value = 4750
$ws.Range("K86").Value = 703
$ws.Range("L86").Value = 4750
$ws.Range("M86").Value = 420
$ws.Range("N86").Value = -6996

# Row 89 (Leve Item ID 12526)
$ws.Range("H89").Value = 2726.5
$ws.Range("I89").Value = 703
$ws.Range("J89").Value = 4750
$ws.Range("K89").Value = 3515
$ws.Range("L89").Value = 23750
$ws.Range("M89").Value = 2101
$ws.Range("N89").Value = -34982

# Row 94 (Leve Item ID 19939)
$ws.Range("H94").Value = 1450.75
$ws.Range("I94").Value = 1436
$ws.Range("J94").Value = 1495
$ws.Range("K94").Value = 1436
$ws.Range("L94").Value = 1495
$ws.Range("M94").Value = -985
$ws.Range("N94").Value = -2397

# Row 134 (Leve Item ID 43998)
$ws.Range("H134").Value = 215386.55
$ws.Range("I134").Value = 1625.4584
$ws.Range("J134").Value = 337535.75
$ws.Range("K134").Value = 4876.3752
$ws.Range("L134").Value = 1012607.25
$ws.Range("M134").Value = -2341.3752
$ws.Range("N134").Value = -1017677.25

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (Leve Item ID 44023)
$ws.Range("H31").Value = 3155.3333
$ws.Range("I31").Value = 2608
$ws.Range("J31").Value = 4250
$ws.Range("K31").Value = 2608
$ws.Range("L31").Value = 4250
$ws.Range("M31").Value = -2313
$ws.Range("N31").Value = -4840

# Row 34 (Leve Item ID 44023)
$ws.Range("H34").Value = 3155.3333
$ws.Range("I34").Value = 2608
$ws.Range("J34").Value = 4250
$ws.Range("K34").Value = 2608
$ws.Range("L34").Value = 4250
$ws.Range("M34").Value = -2406
$ws.Range("N34").Value = -4654

# Row 132 (Leve Item ID 44019)
$ws.Range("H132").Value = 63083.78
$ws.Range("I132").Value = 1784.2354
$ws.Range("J132").Value = 236765.83
$ws.Range("K132").Value = 5352.706200000001
$ws.Range("L132").Value = 710297.49
$ws.Range("M132").Value = -2822.706200000001
$ws.Range("N132").Value = -715357.49

$ws = $wb.Worksheets.Item("CUL")
# Row 12 (Leve Item ID 4854)
$ws.Range("H12").Value = 2
$ws.Range("I12").Value = 2
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 6
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = 167
$ws.Range("N12").Value = ""

# Row 113 (Leve Item ID 27843)
$ws.Range("H113").Value = 800
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 800
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 2400
$ws.Range("M113").Value = ""
$ws.Range("N113").Value = -6740

# Row 118 (Leve Item ID 27872)
$ws.Range("H118").Value = 2924
$ws.Range("I118").Value = 2000
$ws.Range("J118").Value = 3232
$ws.Range("K118").Value = 6000
$ws.Range("L118").Value = 9696
$ws.Range("M118").Value = -4757
$ws.Range("N118").Value = -12182

# Row 121 (Leve Item ID 27878)
$ws.Range("H121").Value = 182819.7
$ws.Range("I121").Value = 233.6
$ws.Range("J121").Value = 443657
$ws.Range("K121").Value = 700.8
$ws.Range("L121").Value = 1330971
$ws.Range("M121").Value = 609.2
$ws.Range("N121").Value = -1333591

$ws = $wb.Worksheets.Item("GSM")
# Row 80 (Leve Item ID 12521)
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = ""
$ws.Range("N80").Value = ""

# Row 83 (Leve Item ID 12521)
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = ""
$ws.Range("N83").Value = ""

# Row 102 (Leve Item ID 36169)
$ws.Range("H102").Value = 1876.4736
$ws.Range("I102").Value = 937.06665
$ws.Range("J102").Value = 5399.25
$ws.Range("K102").Value = 937.06665
$ws.Range("L102").Value = 5399.25
$ws.Range("M102").Value = 684.93335
$ws.Range("N102").Value = -8643.25

# Row 113 (Leve Item ID 27710)
$ws.Range("H113").Value = 1549
$ws.Range("I113").Value = 1549
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1549
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 621
$ws.Range("N113").Value = ""

# Row 122 (Leve Item ID 36182)
$ws.Range("H122").Value = 1593.2
$ws.Range("I122").Value = 1260.5625
$ws.Range("J122").Value = 2923.75
$ws.Range("K122").Value = 3781.6875
$ws.Range("L122").Value = 8771.25
$ws.Range("M122").Value = -1331.6875
$ws.Range("N122").Value = -13671.25

# Row 126 (Leve Item ID 36184)
$ws.Range("H126").Value = 3605.6875
$ws.Range("I126").Value = 3283.8462
$ws.Range("J126").Value = 5000.3335
$ws.Range("K126").Value = 9851.5386
$ws.Range("L126").Value = 15001.0005
$ws.Range("M126").Value = -7381.5386
$ws.Range("N126").Value = -19941.0005

# Row 132 (Leve Item ID 44008)
$ws.Range("H132").Value = 2975.0186
$ws.Range("I132").Value = 1420.2222
$ws.Range("J132").Value = 6084.6113
$ws.Range("K132").Value = 4260.6666
$ws.Range("L132").Value = 18253.8339
$ws.Range("M132").Value = -1730.6666
$ws.Range("N132").Value = -23313.8339

$ws = $wb.Worksheets.Item("LTW")
# Row 132 (Leve Item ID 44058)
$ws.Range("H132").Value = 2116.7678
$ws.Range("I132").Value = 1425.3334
$ws.Range("J132").Value = 4945.364
$ws.Range("K132").Value = 4276.0002
$ws.Range("L132").Value = 14836.092
$ws.Range("M132").Value = -1746.0002
$ws.Range("N132").Value = -19896.092

$ws = $wb.Worksheets.Item("WVR")
# Row 113 (Leve Item ID 27752)
$ws.Range("H113").Value = 441.46155
$ws.Range("I113").Value = 426.9
$ws.Range("J113").Value = 490
$ws.Range("K113").Value = 1280.7
$ws.Range("L113").Value = 1470
$ws.Range("M113").Value = 889.3000000000002
$ws.Range("N113").Value = -5810

